$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the "Femacal de La Calera - Zanahoria"
# data block. It is inserted above the current row 456, shifting the existing
# rows 456:475 down to 457:476 (dimension grows from A1:R475 to A1:R476).
$ws.Rows.Item(456).Insert()

# Populate the newly inserted row 456 with the new record's data.
$ws.Range("A456").Value = 3
$ws.Range("B456").Value = "Femacal de La Calera"
$ws.Range("C456").Value = "Coquimbo"
$ws.Range("D456").Value = 44939
$ws.Range("D456").NumberFormat = $ws.Range("D457").NumberFormat
$ws.Range("E456").Value = 5
$ws.Range("F456").Value = 100114013
$ws.Range("G456").Value = "Zanahoria"
$ws.Range("H456").Value = "Sin especificar"
$ws.Range("I456").Value = "Primera"
$ws.Range("J456").Value = 280
$ws.Range("K456").Value = 11500
$ws.Range("L456").Value = 12000
$ws.Range("M456").Value = 11714
$ws.Range("N456").Value = "$/saco 20 kilos"
$ws.Range("O456").Value = "Provincia de Quillota"
$ws.Range("P456").Value = 586
$ws.Range("Q456").Value = 20
$ws.Range("R456").Value = "Hortaliza"
